# "Add files via upload" - adds a new worksheet "Q1" (with two small cost/
# revenue tables) after "Sheet1", makes it the active sheet, and repositions
# a couple of charts on Sheet1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "Q1" worksheet right after "Sheet1"
# ---------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item(1)
$q1 = $wb.Worksheets.Add($null, $sheet1)
$q1.Name = "Q1"

# ---------------------------------------------------------------------
# 2. First table (rows 4-11) - TFC/TVC/TC/ATC/AVC/MC by unit of output
# ---------------------------------------------------------------------

# Header row (first use of the new shared strings - keeps the same order
# as the source workbook: "Unit (Q)" then "ATC" show up here first)
$q1.Range("A4").Value = "Unit (Q)"
$q1.Range("B4").Value = "TFC"
$q1.Range("C4").Value = "TVC"
$q1.Range("D4").Value = "TC"
$q1.Range("E4").Value = "ATC"
$q1.Range("F4").Value = "AVC"
$q1.Range("G4").Value = "MC"

$q1.Range("A5:G11").HorizontalAlignment = -4131
$q1.Range("E6:F11").NumberFormat = "0.00"

$q1.Range("A5").Value = 0
$q1.Range("B5").Value = 145
$q1.Range("C5").Value = 0
$q1.Range("D5").Formula = "=C5+B5"
$q1.Range("E5").Formula = "=C2"
$q1.Range("F5").Value = 0
$q1.Range("G5").Value = 0

$q1.Range("A6").Value = 1
$q1.Range("B6").Value = 145
$q1.Range("C6").Value = 30
$q1.Range("D6").Formula = "=C6+B6"
$q1.Range("E6").Formula = "=D6/A6"
$q1.Range("F6").Formula = "=C6/A6"
$q1.Range("G6").Formula = "=(D6-D5)/(A6-A5)"

$q1.Range("A7").Value = 2
$q1.Range("B7").Value = 145
$q1.Range("C7").Value = 55
$q1.Range("D7").Formula = "=C7+B7"
$q1.Range("E7").Formula = "=D7/A7"
$q1.Range("F7").Formula = "=C7/A7"
$q1.Range("G7").Formula = "=(D7-D6)/(A7-A6)"

$q1.Range("A8").Value = 3
$q1.Range("B8").Value = 145
$q1.Range("C8").Value = 75
$q1.Range("D8").Formula = "=C8+B8"
$q1.Range("E8").Formula = "=D8/A8"
$q1.Range("F8").Formula = "=C8/A8"
$q1.Range("G8").Formula = "=(D8-D7)/(A8-A7)"

$q1.Range("A9").Value = 4
$q1.Range("B9").Value = 145
$q1.Range("C9").Value = 105
$q1.Range("D9").Formula = "=C9+B9"
$q1.Range("E9").Formula = "=D9/A9"
$q1.Range("F9").Formula = "=C9/A9"
$q1.Range("G9").Formula = "=(D9-D8)/(A9-A8)"

$q1.Range("A10").Value = 5
$q1.Range("B10").Value = 145
$q1.Range("C10").Value = 155
$q1.Range("D10").Formula = "=C10+B10"
$q1.Range("E10").Formula = "=D10/A10"
$q1.Range("F10").Formula = "=C10/A10"
$q1.Range("G10").Formula = "=(D10-D9)/(A10-A9)"

$q1.Range("A11").Value = 6
$q1.Range("B11").Value = 145
$q1.Range("C11").Value = 225
$q1.Range("D11").Formula = "=C11+B11"
$q1.Range("E11").Formula = "=D11/A11"
$q1.Range("F11").Formula = "=C11/A11"
$q1.Range("G11").Formula = "=(D11-D10)/(A11-A10)"

# ---------------------------------------------------------------------
# 3. Second table (rows 18-24) - Price/TR/TC/TFC/TVC/MR/MC by unit
# ---------------------------------------------------------------------

$q1.Range("A18").Value = "Unit (Q)"
$q1.Range("B18").Value = "Price"
$q1.Range("C18").Value = "TR"
$q1.Range("D18").Value = "TC"
$q1.Range("E18").Value = "TFC"
$q1.Range("F18").Value = "TVC"
$q1.Range("G18").Value = "MR"
$q1.Range("H18").Value = "MC"

$q1.Range("A19:H24").HorizontalAlignment = -4131

$q1.Range("A19").Value = 0
$q1.Range("B19").Value = 0
$q1.Range("C19").Value = "-"
$q1.Range("D19").Value = 5
$q1.Range("E19").Value = 5
$q1.Range("F19").Value = 0
$q1.Range("G19").Value = "-"
$q1.Range("H19").Value = "-"

$q1.Range("A20").Value = 1
$q1.Range("B20").Value = 5
$q1.Range("C20").Formula = "=B20*A20"
$q1.Range("D20").Value = 10
$q1.Range("E20").Value = 5
$q1.Range("F20").Formula = "=D20-E20"
$q1.Range("G20").Value = "-"
$q1.Range("H20").Formula = "=(D20-D19)/(A20-A19)"

$q1.Range("A21").Value = 2
$q1.Range("B21").Value = 5
$q1.Range("C21").Formula = "=B21*A21"
$q1.Range("D21").Value = 12
$q1.Range("E21").Value = 5
$q1.Range("F21").Formula = "=D21-E21"
$q1.Range("G21").Formula = "=C21-C20"
$q1.Range("H21").Formula = "=(D21-D20)/(A21-A20)"

$q1.Range("A22").Value = 3
$q1.Range("B22").Value = 5
$q1.Range("C22").Formula = "=B22*A22"
$q1.Range("D22").Value = 15
$q1.Range("E22").Value = 5
$q1.Range("F22").Formula = "=D22-E22"
$q1.Range("G22").Formula = "=C22-C21"
$q1.Range("H22").Formula = "=(D22-D21)/(A22-A21)"

$q1.Range("A23").Value = 4
$q1.Range("B23").Value = 5
$q1.Range("C23").Formula = "=B23*A23"
$q1.Range("D23").Value = 19
$q1.Range("E23").Value = 5
$q1.Range("F23").Formula = "=D23-E23"
$q1.Range("G23").Formula = "=C23-C22"
$q1.Range("H23").Formula = "=(D23-D22)/(A23-A22)"

$q1.Range("A24").Value = 5
$q1.Range("B24").Value = 5
$q1.Range("C24").Formula = "=B24*A24"
$q1.Range("D24").Value = 24
$q1.Range("E24").Value = 5
$q1.Range("F24").Formula = "=D24-E24"
$q1.Range("G24").Formula = "=C24-C23"
$q1.Range("H24").Formula = "=(D24-D23)/(A24-A23)"

# ---------------------------------------------------------------------
# 4. Move / resize the two charts on Sheet1 that shifted position
# ---------------------------------------------------------------------

$chart11 = $sheet1.ChartObjects("Chart 11")
$chart11.Left = 865.589844
$chart11.Top = 1139.25
$chart11.Width = 433.0625
$chart11.Height = 216

$chart6 = $sheet1.ChartObjects("Chart 6")
$chart6.Left = 298.089804
$chart6.Top = 1147.874961
$chart6.Width = 443.5
$chart6.Height = 216

# ---------------------------------------------------------------------
# 5. View state - Q1 becomes the active/selected sheet, Sheet1 keeps a
#    scrolled-down selection.
# ---------------------------------------------------------------------

$sheet1.Range("F63").Select() | Out-Null
$q1.Activate()
$q1.Range("K25").Select() | Out-Null
